$d = $word.ActiveDocument

$replacements = @(
    @("44÷4=", "23÷3="),
    @("15÷5=", "72÷5="),
    @("16÷4=", "77÷8="),
    @("24÷7=", "38÷5="),
    @("70÷3=", "70÷4="),
    @("24÷5=", "17÷2="),
    @("55÷8=", "83÷5="),
    @("89÷8=", "96÷9="),
    @("68÷2=", "35÷5="),
    @("40÷2=", "50÷9="),
    @("57÷5=", "35÷3="),
    @("79÷2=", "91÷7="),
    @("77÷5=", "54÷8="),
    @("21÷3=", "27÷4="),
    @("53÷5=", "83÷3="),
    @("10÷8=", "73÷6="),
    @("19÷4=", "72÷7="),
    @("31÷2=", "57÷2="),
    @("10÷7=", "60÷2="),
    @("64÷9=", "88÷9="),
    @("68÷8=", "18÷6="),
    @("25÷7=", "75÷7="),
    @("85÷7=", "39÷7="),
    @("56÷4=", "16÷7="),
    @("30÷5=", "24÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
